$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.325350284576416
$ws.Range("B1").Value = 3.711971759796143
$ws.Range("C1").Value = 2.04079270362854
$ws.Range("D1").Value = 1.635892152786255
$ws.Range("E1").Value = 1.51152777671814
